$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Insert 11 new rows before the old row 58 (the "Points"/footer block),
#    pushing the existing rows 58-63 down to 69-74. This matches the
#    behaviour visible in the diff: formulas like SUM(D13:D57) keep
#    referencing the original (unshifted) range because the insertion
#    happens entirely below it, while every reference to the footer rows
#    themselves (58->69, 59->70, 60->71, ...) shifts down by 11 rows.
$ws.Rows("58:68").Insert()

# 2) Copy the cell formatting (styles) from the last existing match row
#    (57) down into the newly inserted rows, column by column, so we
#    reuse the same style indices (A/B/E/H/K/N/Q/T/W/Z/AC -> style of
#    A57 etc, C -> style of C57, D/G/J/M/P/S/V/Y/AB -> style of D57 etc.)
#    instead of leaving the blank/duplicated styles that a raw row
#    insert leaves behind.
$allDataCols = @("A","B","C","D","E","G","H","J","K","M","N","P","Q","S","T","V","W","Y","Z","AB","AC")
foreach ($col in $allDataCols) {
    $src = $ws.Range($col + "57")
    $dst = $ws.Range($col + "58:" + $col + "68")
    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = $false

# 3) New match rows (46 through 56), matching the 11 new fixtures added
#    to the schedule.
$matches = @(
    @{Num = 46; Name = "MI vs DC"},
    @{Num = 47; Name = "RR vs CSK"},
    @{Num = 48; Name = "RCB vs PBKS"},
    @{Num = 49; Name = "KKR vs SRH"},
    @{Num = 50; Name = "DC vs CSK"},
    @{Num = 51; Name = "RR vs MI"},
    @{Num = 52; Name = "RCB vs SRH"},
    @{Num = 53; Name = "CSK vs PBKS"},
    @{Num = 54; Name = "KKR vs RR"},
    @{Num = 55; Name = "SRH vs MI"},
    @{Num = 56; Name = "RCB vs DC"}
)

$formulaCols = @("D","G","J","M","P","S","V","Y","AB")
$valueCols   = @("E","H","K","N","Q","T","W","Z","AC")

$row = 58
foreach ($m in $matches) {
    $ws.Range("A$row").Value = $m.Num
    $ws.Range("B$row").Value = 1
    $ws.Range("C$row").Value = $m.Name

    for ($i = 0; $i -lt $formulaCols.Count; $i++) {
        $fcol = $formulaCols[$i]
        $vcol = $valueCols[$i]
        $formula = "=IF(ISERROR(VLOOKUP(RANK($vcol$row, (`$AC$row,`$Z$row,`$W$row,`$T$row,`$Q$row,`$N$row,`$K$row,`$H$row,`$E$row), 0),  `$A`$2:`$B`$10, 2, FALSE)),`"`",VLOOKUP(RANK($vcol$row, (`$AC$row,`$Z$row,`$W$row,`$T$row,`$Q$row,`$N$row,`$K$row,`$H$row,`$E$row), 0),  `$A`$2:`$B`$10, 2, FALSE))"
        $ws.Range("$fcol$row").Formula = $formula
    }

    $row = $row + 1
}

# 4) Restore the view: selection on the totals cell, and keep the frozen
#    pane scrolled near the bottom of the new list.
$ws.Activate()
$ws.Range("A47").Select()
$ws.Range("AD71").Select()

Write-Output "Edit complete"
